$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force text storage for values that would otherwise be auto-coerced
    # to a number by Excel (single-dot decimals, plain integers, etc.),
    # matching the original "Price" column which is stored as text.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 18 and 19 swap: TRON <-> Polkadot (full row content swap with new values)
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D18") "7.25"
$ws.Range("E18").Value = "  -3.60%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D19") "0.115"
$ws.Range("E19").Value = "  -3.64%  "

# Remaining price/volume updates
$ws.Range("D2").Value = "67.754.81"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "3.771.80"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "598.95"
$ws.Range("E5").Value = "  -2.71%  "
Set-TextValue $ws.Range("D6") "174.97"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("D7").Value = "3.770.84"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.10%  "
Set-TextValue $ws.Range("D10") "0.158"
$ws.Range("E10").Value = "  -4.48%  "
Set-TextValue $ws.Range("D11") "6.23"
$ws.Range("E11").Value = "  -5.52%  "
Set-TextValue $ws.Range("D12") "0.466"
$ws.Range("E12").Value = "  -3.70%  "
Set-TextValue $ws.Range("D13") "38.46"
$ws.Range("E13").Value = "  -4.09%  "
Set-TextValue $ws.Range("D14") "0.0000244"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("D15").Value = "4.423.82"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "3.764.03"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "67.805.91"
$ws.Range("E17").Value = "  -2.37%  "
Set-TextValue $ws.Range("D20") "16.78"
$ws.Range("E20").Value = "  +2.81%  "
Set-TextValue $ws.Range("D21") "491.03"
$ws.Range("E21").Value = "  -2.32%  "
Set-TextValue $ws.Range("D22") "9.07"
$ws.Range("E22").Value = "  -2.94%  "
Set-TextValue $ws.Range("D23") "0.748"
$ws.Range("E23").Value = "  +3.35%  "
Set-TextValue $ws.Range("D24") "85.84"
$ws.Range("E24").Value = "  +0.14%  "
Set-TextValue $ws.Range("D25") "0.0000147"
$ws.Range("E25").Value = "  +7.75%  "
Set-TextValue $ws.Range("D26") "2.37"
$ws.Range("E26").Value = "  -7.05%  "
Set-TextValue $ws.Range("D27") "12.35"
$ws.Range("E27").Value = "  -3.59%  "
Set-TextValue $ws.Range("D28") "10.24"
$ws.Range("E28").Value = "  -4.21%  "
Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.37%  "
Set-TextValue $ws.Range("D31") "2.43"
$ws.Range("E31").Value = "  -3.48%  "
Set-TextValue $ws.Range("D32") "33.00"
$ws.Range("E32").Value = "  +8.58%  "
Set-TextValue $ws.Range("D33") "7.73"
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("E34").Value = "  -3.65%  "
Set-TextValue $ws.Range("D35") "0.999"
$ws.Range("E35").Value = "  -0.08%  "
Set-TextValue $ws.Range("D36") "1.01"
$ws.Range("E36").Value = "  -3.40%  "
Set-TextValue $ws.Range("D37") "5.82"
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("E38").Value = "  -2.45%  "
Set-TextValue $ws.Range("D39") "0.327"
$ws.Range("E39").Value = "  -4.60%  "
Set-TextValue $ws.Range("D40") "452.73"
$ws.Range("E40").Value = "  +1.88%  "
Set-TextValue $ws.Range("D41") "49.21"
$ws.Range("E41").Value = "  -0.89%  "
Set-TextValue $ws.Range("D42") "2.01"
$ws.Range("E42").Value = "  -3.04%  "
Set-TextValue $ws.Range("D43") "2.89"
$ws.Range("E43").Value = "  -2.98%  "
Set-TextValue $ws.Range("D44") "8.40"
$ws.Range("E44").Value = "  -1.84%  "
Set-TextValue $ws.Range("D45") "40.83"
$ws.Range("E45").Value = "  -8.54%  "
$ws.Range("D46").Value = "2.833.40"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("E47").Value = "  -0.04%  "
Set-TextValue $ws.Range("D48") "138.56"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("E49").Value = "  -1.92%  "
Set-TextValue $ws.Range("D50") "26.07"
$ws.Range("E50").Value = "  -5.24%  "
Set-TextValue $ws.Range("D51") "23.50"
$ws.Range("E51").Value = "  +5.17%  "
